$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 2 de Abril de 2020 a las 09:25'

$ws.Range('B8').Value = 77981
$ws.Range('C8').Value = 0
$ws.Range('D8').Value = 19175
$ws.Range('E8').Value = 57875
$ws.Range('F8').Value = 3408
$ws.Range('G8').Value = 0
$ws.Range('H8').Value = 931

$ws.Range('B16').Value = 10809
$ws.Range('C16').Value = 98
$ws.Range('D16').Value = 1436
$ws.Range('E16').Value = 9227
$ws.Range('F16').Value = 215
$ws.Range('G16').Value = 0
$ws.Range('H16').Value = 146

$ws.Range('B21').Value = 6211
$ws.Range('C21').Value = 119
$ws.Range('D21').Value = 289
$ws.Range('E21').Value = 5892
$ws.Range('F21').Value = 107
$ws.Range('G21').Value = 4
$ws.Range('H21').Value = 30

$ws.Range('B22').Value = 5137
$ws.Range('C22').Value = 89
$ws.Range('D22').Value = 345
$ws.Range('E22').Value = 4767
$ws.Range('F22').Value = 50
$ws.Range('G22').Value = 2
$ws.Range('H22').Value = 25

$ws.Range('B25').Value = 3604
$ws.Range('C25').Value = 15
$ws.Range('D25').Value = 61
$ws.Range('E25').Value = 3503
$ws.Range('F25').Value = 72
$ws.Range('G25').Value = 1
$ws.Range('H25').Value = 40

$ws.Range('B33').Value = 2460
$ws.Range('C33').Value = 0
$ws.Range('D33').Value = 252
$ws.Range('E33').Value = 2114
$ws.Range('F33').Value = 57
$ws.Range('G33').Value = 2
$ws.Range('H33').Value = 94

$ws.Range('B42').Value = 1518
$ws.Range('C42').Value = 72
$ws.Range('D42').Value = 300
$ws.Range('E42').Value = 1201
$ws.Range('F42').Value = 62
$ws.Range('G42').Value = 0
$ws.Range('H42').Value = 17

$ws.Range('A59').Value = 'Ucrania'
$ws.Range('B59').Value = 804
$ws.Range('C59').Value = 10
$ws.Range('D59').Value = 13
$ws.Range('E59').Value = 771
$ws.Range('F59').Value = 0
$ws.Range('G59').Value = 0
$ws.Range('H59').Value = 20

$ws.Range('A60').Value = 'Nueva Zelanda'
$ws.Range('B60').Value = 797
$ws.Range('C60').Value = 89
$ws.Range('D60').Value = 92
$ws.Range('E60').Value = 704
$ws.Range('F60').Value = 2
$ws.Range('G60').Value = 0
$ws.Range('H60').Value = 1

$ws.Range('A66').Value = 'Armenia'
$ws.Range('B66').Value = 663
$ws.Range('C66').Value = 92
$ws.Range('D66').Value = 33
$ws.Range('E66').Value = 626
$ws.Range('F66').Value = 30
$ws.Range('G66').Value = 0
$ws.Range('H66').Value = 4

$ws.Range('A67').Value = 'Marruecos'
$ws.Range('B67').Value = 654
$ws.Range('C67').Value = 0
$ws.Range('D67').Value = 29
$ws.Range('E67').Value = 586
$ws.Range('F67').Value = 1
$ws.Range('G67').Value = 0
$ws.Range('H67').Value = 39

$ws.Range('A68').Value = 'Lituania'
$ws.Range('B68').Value = 649
$ws.Range('C68').Value = 68
$ws.Range('D68').Value = 7
$ws.Range('E68').Value = 634
$ws.Range('F68').Value = 11
$ws.Range('G68').Value = 0
$ws.Range('H68').Value = 8

$ws.Range('A69').Value = 'Hungria'
$ws.Range('B69').Value = 585
$ws.Range('C69').Value = 60
$ws.Range('D69').Value = 42
$ws.Range('E69').Value = 522
$ws.Range('F69').Value = 17
$ws.Range('G69').Value = 1
$ws.Range('H69').Value = 21

$ws.Range('A71').Value = 'Bosnia y Herzegovina'
$ws.Range('B71').Value = 490
$ws.Range('C71').Value = 31
$ws.Range('D71').Value = 19
$ws.Range('E71').Value = 458
$ws.Range('F71').Value = 1
$ws.Range('G71').Value = 0
$ws.Range('H71').Value = 13

$ws.Range('A72').Value = 'Libano'
$ws.Range('B72').Value = 479
$ws.Range('C72').Value = 0
$ws.Range('D72').Value = 43
$ws.Range('E72').Value = 422
$ws.Range('F72').Value = 5
$ws.Range('G72').Value = 0
$ws.Range('H72').Value = 14

$ws.Range('B77').Value = 402
$ws.Range('C77').Value = 22
$ws.Range('D77').Value = 27
$ws.Range('E77').Value = 372
$ws.Range('F77').Value = 6
$ws.Range('G77').Value = 0
$ws.Range('H77').Value = 3

$ws.Range('A94').Value = 'Oman'
$ws.Range('B94').Value = 231
$ws.Range('C94').Value = 21
$ws.Range('D94').Value = 41
$ws.Range('E94').Value = 189
$ws.Range('F94').Value = 3
$ws.Range('G94').Value = 0
$ws.Range('H94').Value = 1

$ws.Range('A95').Value = 'Vietnam'
$ws.Range('B95').Value = 222
$ws.Range('C95').Value = 4
$ws.Range('D95').Value = 64
$ws.Range('E95').Value = 158
$ws.Range('F95').Value = 3
$ws.Range('G95').Value = 0
$ws.Range('H95').Value = 0

$ws.Range('A96').Value = 'Honduras'
$ws.Range('B96').Value = 219
$ws.Range('C96').Value = 47
$ws.Range('D96').Value = 3
$ws.Range('E96').Value = 202
$ws.Range('F96').Value = 4
$ws.Range('G96').Value = 4
$ws.Range('H96').Value = 14

$ws.Range('A97').Value = 'Cuba'
$ws.Range('B97').Value = 212
$ws.Range('C97').Value = 0
$ws.Range('D97').Value = 12
$ws.Range('E97').Value = 194
$ws.Range('F97').Value = 3
$ws.Range('G97').Value = 0
$ws.Range('H97').Value = 6

$ws.Range('A109').Value = 'Montenegro'
$ws.Range('B109').Value = 140
$ws.Range('C109').Value = 17
$ws.Range('D109').Value = 0
$ws.Range('E109').Value = 138
$ws.Range('F109').Value = 4
$ws.Range('G109').Value = 0
$ws.Range('H109').Value = 2

$ws.Range('A110').Value = 'Martinica'
$ws.Range('B110').Value = 135
$ws.Range('C110').Value = 0
$ws.Range('D110').Value = 27
$ws.Range('E110').Value = 105
$ws.Range('F110').Value = 16
$ws.Range('G110').Value = 0
$ws.Range('H110').Value = 3

$ws.Range('A111').Value = 'Estado de Palestina'
$ws.Range('B111').Value = 134
$ws.Range('C111').Value = 0
$ws.Range('D111').Value = 18
$ws.Range('E111').Value = 115
$ws.Range('F111').Value = 0
$ws.Range('G111').Value = 0
$ws.Range('H111').Value = 1

$ws.Range('A112').Value = 'Brunei'
$ws.Range('B112').Value = 131
$ws.Range('C112').Value = 0
$ws.Range('D112').Value = 52
$ws.Range('E112').Value = 78
$ws.Range('F112').Value = 3
$ws.Range('G112').Value = 0
$ws.Range('H112').Value = 1

$ws.Range('A113').Value = 'Georgia'
$ws.Range('B113').Value = 130
$ws.Range('C113').Value = 13
$ws.Range('D113').Value = 23
$ws.Range('E113').Value = 107
$ws.Range('F113').Value = 6
$ws.Range('G113').Value = 0
$ws.Range('H113').Value = 0

$ws.Range('A114').Value = 'Guadalupe'
$ws.Range('B114').Value = 125
$ws.Range('C114').Value = 0
$ws.Range('D114').Value = 24
$ws.Range('E114').Value = 95
$ws.Range('F114').Value = 14
$ws.Range('G114').Value = 0
$ws.Range('H114').Value = 6
